# Update "弹幕数" (barrage/view-count style) figures in column F for a handful
# of rows on sheet "展览" (Exhibition), "本地生活" (Local life), and
# "全部类型" (All types) to reflect freshly generated site stats.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- 展览 (sheet1) ---
$wsExhibit.Range("F2").Value  = 601
$wsExhibit.Range("F6").Value  = 406
$wsExhibit.Range("F8").Value  = 174
$wsExhibit.Range("F10").Value = 241
$wsExhibit.Range("F11").Value = 6226
$wsExhibit.Range("F13").Value = 56
$wsExhibit.Range("F14").Value = 523
$wsExhibit.Range("F27").Value = 1888

# --- 本地生活 (sheet3) ---
$wsLocal.Range("F2").Value = 288

# --- 全部类型 (sheet4) ---
$wsAll.Range("F2").Value  = 288
$wsAll.Range("F3").Value  = 601
$wsAll.Range("F8").Value  = 406
$wsAll.Range("F10").Value = 174
$wsAll.Range("F12").Value = 241
$wsAll.Range("F13").Value = 6226
$wsAll.Range("F15").Value = 56
$wsAll.Range("F17").Value = 523
$wsAll.Range("F37").Value = 1888
